$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 24.84807066666667
$ws.Range("H2").Value = 74.544212
$ws.Range("I2").Value = 0.1383130136760648
$ws.Range("J2").Value = 0.1383130136760648
$ws.Range("M2").Value = 13.929953
$ws.Range("N2").Value = 41.789859
$ws.Range("O2").Value = 0.09674275490334808
$ws.Range("P2").Value = 0.09674275490334808
$ws.Range("Q2").Value = 346.1324565273454
$ws.Range("R2").Value = 3115.192108746108
$ws.Range("S2").Value = 0.01338078198200697
$ws.Range("T2").Value = 0.01338078198200697
$ws.Range("G3").Value = 24.84807066666667
$ws.Range("H3").Value = 74.544212
$ws.Range("I3").Value = 0.1383130136760648
$ws.Range("J3").Value = 0.1383130136760648
$ws.Range("M3").Value = 81.07766966666667
$ws.Range("O3").Value = 0.5630799418129374
$ws.Range("P3").Value = 0.5630799418129373
$ws.Range("Q3").Value = 2014.62366536599
$ws.Range("R3").Value = 18131.61298829391
$ws.Range("S3").Value = 0.07788128369269059
$ws.Range("T3").Value = 0.07788128369269058
$ws.Range("G4").Value = 24.84807066666667
$ws.Range("H4").Value = 74.544212
$ws.Range("I4").Value = 0.1383130136760648
$ws.Range("J4").Value = 0.1383130136760648
$ws.Range("M4").Value = 48.98200233333333
$ws.Range("N4").Value = 146.946007
$ws.Range("O4").Value = 0.3401773032837146
$ws.Range("P4").Value = 0.3401773032837146
$ws.Range("Q4").Value = 1217.108255373498
$ws.Range("R4").Value = 10953.97429836149
$ws.Range("S4").Value = 0.04705094800136728
$ws.Range("T4").Value = 0.04705094800136728
$ws.Range("I5").Value = 0.5888817337812031
$ws.Range("J5").Value = 0.5888817337812031
$ws.Range("M5").Value = 13.929953
$ws.Range("N5").Value = 41.789859
$ws.Range("O5").Value = 0.09674275490334808
$ws.Range("P5").Value = 0.09674275490334808
$ws.Range("Q5").Value = 1473.69416442007
$ws.Range("R5").Value = 13263.24747978063
$ws.Range("S5").Value = 0.0569700412382536
$ws.Range("T5").Value = 0.0569700412382536
$ws.Range("I6").Value = 0.5888817337812031
$ws.Range("J6").Value = 0.5888817337812031
$ws.Range("M6").Value = 81.07766966666667
$ws.Range("O6").Value = 0.5630799418129374
$ws.Range("P6").Value = 0.5630799418129373
$ws.Range("Q6").Value = 8577.465311802902
$ws.Range("R6").Value = 77197.18780622611
$ws.Range("S6").Value = 0.3315874923922216
$ws.Range("T6").Value = 0.3315874923922215
$ws.Range("I7").Value = 0.5888817337812031
$ws.Range("J7").Value = 0.5888817337812031
$ws.Range("M7").Value = 48.98200233333333
$ws.Range("N7").Value = 146.946007
$ws.Range("O7").Value = 0.3401773032837146
$ws.Range("P7").Value = 0.3401773032837146
$ws.Range("Q7").Value = 5181.962279430776
$ws.Range("R7").Value = 46637.66051487698
$ws.Range("S7").Value = 0.200324200150728
$ws.Range("T7").Value = 0.200324200150728
$ws.Range("G8").Value = 49.00973533333333
$ws.Range("H8").Value = 147.029206
$ws.Range("I8").Value = 0.272805252542732
$ws.Range("J8").Value = 0.272805252542732
$ws.Range("M8").Value = 13.929953
$ws.Range("N8").Value = 41.789859
$ws.Range("O8").Value = 0.09674275490334808
$ws.Range("P8").Value = 0.09674275490334808
$ws.Range("Q8").Value = 682.7033097357726
$ws.Range("R8").Value = 6144.329787621953
$ws.Range("S8").Value = 0.02639193168308749
$ws.Range("T8").Value = 0.02639193168308749
$ws.Range("G9").Value = 49.00973533333333
$ws.Range("H9").Value = 147.029206
$ws.Range("I9").Value = 0.272805252542732
$ws.Range("J9").Value = 0.272805252542732
$ws.Range("M9").Value = 81.07766966666667
$ws.Range("O9").Value = 0.5630799418129374
$ws.Range("P9").Value = 0.5630799418129373
$ws.Range("Q9").Value = 3973.595131806761
$ws.Range("R9").Value = 35762.35618626085
$ws.Range("S9").Value = 0.1536111657280252
$ws.Range("T9").Value = 0.1536111657280252
$ws.Range("G10").Value = 49.00973533333333
$ws.Range("H10").Value = 147.029206
$ws.Range("I10").Value = 0.272805252542732
$ws.Range("J10").Value = 0.272805252542732
$ws.Range("M10").Value = 48.98200233333333
$ws.Range("N10").Value = 146.946007
$ws.Range("O10").Value = 0.3401773032837146
$ws.Range("P10").Value = 0.3401773032837146
$ws.Range("Q10").Value = 2400.594970453382
$ws.Range("R10").Value = 21605.35473408044
$ws.Range("S10").Value = 0.0928021551316193
$ws.Range("T10").Value = 0.0928021551316193